$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'93.441.54"
$ws.Range("E2").Value = "  -4.13%  "

$ws.Range("D3").Value = "'3.398.94"
$ws.Range("E3").Value = "  +1.67%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'233.94"
$ws.Range("E5").Value = "  -7.08%  "

$ws.Range("D6").Value = "'632.70"
$ws.Range("E6").Value = "  -3.89%  "

$ws.Range("D7").Value = "'1.41"
$ws.Range("E7").Value = "  -1.77%  "

$ws.Range("D8").Value = "'0.389"
$ws.Range("E8").Value = "  -8.16%  "

$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("D10").Value = "'0.948"
$ws.Range("E10").Value = "  -6.89%  "

$ws.Range("D11").Value = "'3.394.20"
$ws.Range("E11").Value = "  +1.58%  "

$ws.Range("D12").Value = "'41.41"
$ws.Range("E12").Value = "  +0.11%  "

$ws.Range("D13").Value = "'0.196"
$ws.Range("E13").Value = "  -5.49%  "

$ws.Range("D14").Value = "'6.06"
$ws.Range("E14").Value = "  -0.74%  "

$ws.Range("D15").Value = "'93.419.88"
$ws.Range("E15").Value = "  -3.93%  "

$ws.Range("D16").Value = "'4.036.21"
$ws.Range("E16").Value = "  +1.46%  "

$ws.Range("D17").Value = "'0.0000247"
$ws.Range("E17").Value = "  -2.84%  "

$ws.Range("D18").Value = "'8.25"
$ws.Range("E18").Value = "  -6.31%  "

$ws.Range("D19").Value = "'3.404.54"
$ws.Range("E19").Value = "  +2.61%  "

$ws.Range("D20").Value = "'17.37"
$ws.Range("E20").Value = "  -2.40%  "

$ws.Range("D21").Value = "'11.10"
$ws.Range("E21").Value = "  +3.45%  "

$ws.Range("D22").Value = "'0.491"
$ws.Range("E22").Value = "  -11.25%  "

$ws.Range("D23").Value = "'490.42"
$ws.Range("E23").Value = "  -4.10%  "

$ws.Range("D24").Value = "'3.16"
$ws.Range("E24").Value = "  -5.64%  "

$ws.Range("D25").Value = "'6.45"
$ws.Range("E25").Value = "  -2.39%  "

$ws.Range("D26").Value = "'0.0000186"
$ws.Range("E26").Value = "  -7.20%  "

$ws.Range("D27").Value = "'90.42"
$ws.Range("E27").Value = "  -7.24%  "

$ws.Range("D28").Value = "'3.579.37"
$ws.Range("E28").Value = "  +1.29%  "

$ws.Range("D29").Value = "'11.83"
$ws.Range("E29").Value = "  -3.38%  "

$ws.Range("D30").Value = "'11.49"
$ws.Range("E30").Value = "  -1.54%  "

$ws.Range("D31").Value = "'0.994"
$ws.Range("E31").Value = "  -0.64%  "

$ws.Range("D32").Value = "'2.69"
$ws.Range("E32").Value = "  +6.15%  "

$ws.Range("D33").Value = "'0.134"
$ws.Range("E33").Value = "  -9.20%  "

$ws.Range("D34").Value = "'0.179"
$ws.Range("E34").Value = "  -6.21%  "

$ws.Range("D35").Value = "'0.996"
$ws.Range("E35").Value = "  -0.27%  "

$ws.Range("D36").Value = "'29.68"
$ws.Range("E36").Value = "  +4.00%  "

$ws.Range("D37").Value = "'0.549"
$ws.Range("E37").Value = "  -2.18%  "

$ws.Range("D38").Value = "'533.63"
$ws.Range("E38").Value = "  +4.18%  "

$ws.Range("D39").Value = "'7.49"
$ws.Range("E39").Value = "  -5.35%  "

$ws.Range("D40").Value = "'1.42"
$ws.Range("E40").Value = "  -4.49%  "

$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.149"
$ws.Range("E42").Value = "  -2.02%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'0.913"
$ws.Range("E43").Value = "  +8.02%  "

$ws.Range("D44").Value = "'24.03"
$ws.Range("E44").Value = "  -1.51%  "

$ws.Range("D45").Value = "'1.65"
$ws.Range("E45").Value = "  -2.31%  "

$ws.Range("D46").Value = "'0.0406"
$ws.Range("E46").Value = "  -7.00%  "

$ws.Range("D47").Value = "'5.43"
$ws.Range("E47").Value = "  -4.50%  "

$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'53.17"
$ws.Range("E48").Value = "  -2.13%  "

$ws.Range("B49").Value = "MantraDAO"
$ws.Range("C49").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D49").Value = "'3.41"
$ws.Range("E49").Value = "  -6.19%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'2.11"
$ws.Range("E50").Value = "  +5.14%  "

$ws.Range("D51").Value = "'3.11"
$ws.Range("E51").Value = "  -1.28%  "
